$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B23 text: self-reference corrected from "e021" to "e022" (Choose Operations)
$ws.Range("B23").Value = "<Bold>e022 Choose Operations</Bold> `n<InlineUIContainer><Button Content='r4.54' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nChoose one of following options. To see the options, select the &apos;e###&apos; button. To choose the option, select the other buttons. Each option uses up time per the `n<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nTable.`n<LineBreak/><LineBreak/>`n<InlineUIContainer><Button Content='e020' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`n Perform an additional Enemy <InlineUIContainer><Button Content='Strength Check' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>`n<InlineUIContainer><Button Content='e023' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Call for Artillery `n<InlineUIContainer><Button Content='Support' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>`n<InlineUIContainer><Button Content='e025' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Call for Air `n<InlineUIContainer><Button Content='Strike' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>`n<InlineUIContainer><Button Content='e027' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Attempt to `n<InlineUIContainer><Button Content='Resupply' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>`n<InlineUIContainer><Button Content='e028' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Choose adjacent area to `n<InlineUIContainer><Button Content='Enter' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>"

# Row 23 height shrinks by one wrapped line now that the trailing blank line was removed
$ws.Range("B23").EntireRow.RowHeight = 255

# Update the visible selection to match the saved view state
$ws.Activate()
$ws.Range("B24").Select()
